$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$timestamps = @(
    "Thu Nov 07 16:10:46 EST 2024",
    "Thu Nov 07 16:10:56 EST 2024",
    "Thu Nov 07 16:11:07 EST 2024",
    "Thu Nov 07 16:11:19 EST 2024",
    "Thu Nov 07 16:11:30 EST 2024",
    "Thu Nov 07 16:11:42 EST 2024",
    "Thu Nov 07 16:11:53 EST 2024",
    "Thu Nov 07 16:12:04 EST 2024",
    "Thu Nov 07 16:12:15 EST 2024",
    "Thu Nov 07 16:12:26 EST 2024",
    "Thu Nov 07 16:12:37 EST 2024",
    "Thu Nov 07 16:12:48 EST 2024",
    "Thu Nov 07 16:12:59 EST 2024",
    "Thu Nov 07 16:13:10 EST 2024",
    "Thu Nov 07 16:13:21 EST 2024",
    "Thu Nov 07 16:13:35 EST 2024",
    "Thu Nov 07 16:13:47 EST 2024",
    "Thu Nov 07 16:14:00 EST 2024",
    "Thu Nov 07 16:14:13 EST 2024",
    "Thu Nov 07 16:14:25 EST 2024",
    "Thu Nov 07 16:14:37 EST 2024",
    "Thu Nov 07 16:14:50 EST 2024",
    "Thu Nov 07 16:15:02 EST 2024",
    "Thu Nov 07 16:15:13 EST 2024",
    "Thu Nov 07 16:15:24 EST 2024",
    "Thu Nov 07 16:15:35 EST 2024",
    "Thu Nov 07 16:15:46 EST 2024",
    "Thu Nov 07 16:15:57 EST 2024",
    "Thu Nov 07 16:16:08 EST 2024",
    "Thu Nov 07 16:16:19 EST 2024",
    "Thu Nov 07 16:16:30 EST 2024",
    "Thu Nov 07 16:16:41 EST 2024",
    "Thu Nov 07 16:16:52 EST 2024",
    "Thu Nov 07 16:17:04 EST 2024",
    "Thu Nov 07 16:17:15 EST 2024",
    "Thu Nov 07 16:17:26 EST 2024",
    "Thu Nov 07 16:17:38 EST 2024",
    "Thu Nov 07 16:17:50 EST 2024",
    "Thu Nov 07 16:18:02 EST 2024",
    "Thu Nov 07 16:18:13 EST 2024",
    "Thu Nov 07 16:18:24 EST 2024",
    "Thu Nov 07 16:18:35 EST 2024",
    "Thu Nov 07 16:18:47 EST 2024",
    "Thu Nov 07 16:19:00 EST 2024",
    "Thu Nov 07 16:19:11 EST 2024",
    "Thu Nov 07 16:19:22 EST 2024",
    "Thu Nov 07 16:19:33 EST 2024"
)

for ($i = 0; $i -lt $timestamps.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $timestamps[$i]
}

